$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$row = $table.Rows.Item(6)

$row.Cells.Item(1).Range.Text = "5"
$row.Cells.Item(2).Range.Text = "House001"
$row.Cells.Item(3).Range.Text = "HOUSE"
$row.Cells.Item(4).Range.Text = "LightBlazeMC - Turok"
$row.Cells.Item(5).Range.Text = "Owned by me"
